$d = $word.ActiveDocument

# 1) Update the arraignment date reference and the fines/costs due date
#    from December 22, 2021 to December 23, 2021 (both occurrences).
$d.Content.Find.Execute("December 22, 2021", $false, $false, $false, $false, $false,
                         $true, 1, $false, "December 23, 2021", 2)

# 2) Change the Plea row entries from "Guilty" to "No Contest" for the
#    four charges, leaving the Finding row ("Guilty") untouched.
for ($col = 2; $col -le 5; $col++) {
    $table = $d.Tables.Item(1)
    $cell = $table.Cell(4, $col)
    $cell.Range.Text = "No Contest"
}
